$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell as literal TEXT (matching the rest of the
# sheet, which stores everything - including numeric-looking values like
# "120" - as text) instead of letting Excel auto-infer a number type for
# digit-only strings. A plain "$ws.Range(...).Value = '123'" assignment
# would be auto-coerced to a number by Excel, so instead we build the text
# in a scratch cell via a formula (which always yields a text result),
# freeze it to a static value with PasteSpecial (values only), and copy
# that static text into the destination cell. This preserves the text type
# without applying any text number-format / quote-prefix style to the cell.
function Set-TextValue {
    param($cellRef, [string]$text)

    $helper = $ws.Range("Z100")
    $helper.Formula = '=""&"' + $text + '"'
    $helper.Copy()
    $helper.PasteSpecial(-4163)  # xlPasteValues - freeze formula to static text
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues - copy as text
    $helper.ClearContents()
}

# Row 4 ("asd asd asd asd asd 120 ...") gains values across G:M
Set-TextValue "G4" "123"
Set-TextValue "H4" "123"
Set-TextValue "I4" "123"
Set-TextValue "J4" "123"
Set-TextValue "K4" "1231"
Set-TextValue "L4" "123"
Set-TextValue "M4" "123"

# Row 5's previously-blank I5/K5 cells are removed entirely
$ws.Range("I5").ClearContents()
$ws.Range("K5").ClearContents()
